# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets, matching the regenerated data
# output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F holds the interested-count value.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 96
$ws1.Range("F3").Value = 821
$ws1.Range("F6").Value = 135
$ws1.Range("F8").Value = 4833
$ws1.Range("F10").Value = 5168
$ws1.Range("F12").Value = 1292
$ws1.Range("F13").Value = 94

# Sheet "全部类型" (all types) mirrors the same rows (with one extra
# row inserted at position 7), so the row offsets differ slightly.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 96
$ws4.Range("F3").Value = 821
$ws4.Range("F6").Value = 135
$ws4.Range("F9").Value = 4833
$ws4.Range("F11").Value = 5168
$ws4.Range("F13").Value = 1292
$ws4.Range("F14").Value = 94
